$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts existing D..H to E..I)
$ws.Columns.Item(4).Insert()

# Header for the new column
$ws.Range("D1").Value = "Tavg"

# Formula for the new column: average of Tn (B) and Tx (C)
$ws.Range("D2").Formula = "=(B2+C2)/2"
$ws.Range("D3:D41").Formula = "=(B3+C3)/2"

# Update the view to match target state
$ws.Range("D2:D41").Select()
$excel.ActiveWindow.ScrollRow = 26
